$d = $word.ActiveDocument

$oldText = "W.S. Brown and M. Grossi (2016) Pre- and Post-Mission-6 Glider CTD Comparison Measurements: 11 June and 22 July 2015, University of Massachusetts Dartmouth School for Marine Science and Technology Technical Report SMAST-16-0501."
$newText = "W.S. Brown and M. Grossi (2015) Pre- and Post-Mission Glider CTD Comparison Measurements: 19 June 2014 and 6 February 2015, University of Massachusetts Dartmouth School for Marine Science and Technology Technical Report SMAST-15-06-01."

# There are two identical paragraphs with this citation; only the second
# (last) occurrence - the one immediately preceding bookmarkEnd id=32 - is
# updated by this change. Find it by locating every paragraph whose text
# matches the old citation and operating on the final match.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    $pText = $p.Range.Text
    $pText = $pText.TrimEnd([char]13, [char]7)
    if ($pText -eq $oldText) {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $r = $targetPara.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}
